$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC sheet - row 38
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2616.65
$ws.Range("I38").Value = 121.181816
$ws.Range("J38").Value = 5666.6665
$ws.Range("K38").Value = 363.545448
$ws.Range("L38").Value = 16999.9995
$ws.Range("M38").Value = 8.454552000000035
$ws.Range("N38").Value = -17743.9995

# ---------------------------------------------------------------------------
# ARM sheet - rows 118, 132
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H132").Value = 133066.88
$ws.Range("I132").Value = 144158.55
$ws.Range("J132").Value = 3664
$ws.Range("K132").Value = 432475.65
$ws.Range("L132").Value = 10992
$ws.Range("M132").Value = -429945.65
$ws.Range("N132").Value = -16052

# ---------------------------------------------------------------------------
# CRP sheet - rows 16, 32, 113, 125
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 1842.75
$ws.Range("I16").Value = 1662.5
$ws.Range("K16").Value = 1662.5
$ws.Range("M16").Value = -1375.5

$ws.Range("H32").Value = 6100
$ws.Range("I32").Value = 6100
$ws.Range("K32").Value = 6100
$ws.Range("M32").Value = -5784

$ws.Range("H113").Value = 1842.75
$ws.Range("I113").Value = 1662.5
$ws.Range("K113").Value = 1662.5
$ws.Range("M113").Value = 507.5

$ws.Range("H125").Value = 34996.668
$ws.Range("J125").Value = 34996.668
$ws.Range("L125").Value = 34996.668
$ws.Range("N125").Value = -39916.668

# ---------------------------------------------------------------------------
# CUL sheet - rows 41, 44, 131
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H41").Value = 825
$ws.Range("I41").Value = 500
$ws.Range("J41").Value = 933.3333
$ws.Range("K41").Value = 1500
$ws.Range("L41").Value = 2799.9999
$ws.Range("M41").Value = -1162
$ws.Range("N41").Value = -3475.9999

$ws.Range("H44").Value = 699.1087
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 699.1087
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 2097.3261
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -2893.3261

$ws.Range("H131").Value = 951.6957
$ws.Range("J131").Value = 982.0158699999999
$ws.Range("L131").Value = 2946.04761
$ws.Range("N131").Value = -13026.04761

# ---------------------------------------------------------------------------
# LTW sheet - rows 114, 120, 124-141
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H114").Value = 1400
$ws.Range("J114").Value = 1400
$ws.Range("L114").Value = 1400
$ws.Range("N114").Value = -10078

$ws.Range("H120").Value = 31000
$ws.Range("J120").Value = 31000
$ws.Range("L120").Value = 31000
$ws.Range("N120").Value = -40676

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0

$ws.Range("H125").Value = 32398.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 32398.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 32398.5
$ws.Range("N125").Value = -42238.5

$ws.Range("H127").Value = 46330
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 46330
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 46330
$ws.Range("N127").Value = -56250

$ws.Range("H128").Value = 45999.668
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 45999.668
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 45999.668
$ws.Range("N128").Value = -55959.668

$ws.Range("H129").Value = 44940
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 44940
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 44940
$ws.Range("N129").Value = -54940

$ws.Range("H130").Value = 58681.8
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 58681.8
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 58681.8
$ws.Range("N130").Value = -68721.8

$ws.Range("H131").Value = 25368
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 25368
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 25368
$ws.Range("N131").Value = -35448

$ws.Range("H132").Value = 2510.0344
$ws.Range("I132").Value = 2252.2354
$ws.Range("J132").Value = 2875.25
$ws.Range("K132").Value = 6756.706200000001
$ws.Range("L132").Value = 8625.75
$ws.Range("M132").Value = -4226.706200000001
$ws.Range("N132").Value = -13685.75

$ws.Range("H133").Value = 54906.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 54906.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 54906.5
$ws.Range("N133").Value = -59966.5

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0

$ws.Range("H135").Value = 100429
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 100429
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 100429
$ws.Range("N135").Value = -110569

$ws.Range("H136").Value = 1594.3269
$ws.Range("I136").Value = 1365
$ws.Range("J136").Value = 2216.7856
$ws.Range("K136").Value = 4095
$ws.Range("L136").Value = 6650.3568
$ws.Range("M136").Value = -1545
$ws.Range("N136").Value = -11750.3568

$ws.Range("H137").Value = 53429
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 53429
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 53429
$ws.Range("N137").Value = -63629

$ws.Range("H138").Value = 41325.777
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 41325.777
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 41325.777
$ws.Range("N138").Value = -51605.777

$ws.Range("H139").Value = 39370
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 39370
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 39370
$ws.Range("N139").Value = -49650

$ws.Range("H140").Value = 45404.668
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 45404.668
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 45404.668
$ws.Range("N140").Value = -55764.668

$ws.Range("H141").Value = 73876.375
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 73876.375
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 73876.375
$ws.Range("N141").Value = -84236.375

# ---------------------------------------------------------------------------
# WVR sheet - row 111
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H111").Value = 36888.6
$ws.Range("J111").Value = 36888.6
$ws.Range("L111").Value = 36888.6
$ws.Range("N111").Value = -45068.6
